$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.103903333333333
$ws.Range("H2").Value = 3.31171
$ws.Range("I2").Value = 0.02393122995918198
$ws.Range("J2").Value = 0.02393122995918198
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 13.10121233333333
$ws.Range("N2").Value = 39.303637
$ws.Range("O2").Value = 0.1081423012186565
$ws.Range("P2").Value = 0.1081423012186565
$ws.Range("Q2").Value = 14.46247196547445
$ws.Range("R2").Value = 130.16224768927
$ws.Range("S2").Value = 0.002587978278778793
$ws.Range("T2").Value = 0.002587978278778794
$ws.Range("G3").Value = 1.103903333333333
$ws.Range("H3").Value = 3.31171
$ws.Range("I3").Value = 0.02393122995918198
$ws.Range("J3").Value = 0.02393122995918198
$ws.Range("O3").Value = 0.5751439322003361
$ws.Range("P3").Value = 0.5751439322003362
$ws.Range("Q3").Value = 76.91719985449222
$ws.Range("R3").Value = 692.2547986904299
$ws.Range("S3").Value = 0.01376390170111441
$ws.Range("T3").Value = 0.01376390170111442
$ws.Range("G4").Value = 1.103903333333333
$ws.Range("H4").Value = 3.31171
$ws.Range("I4").Value = 0.02393122995918198
$ws.Range("J4").Value = 0.02393122995918198
$ws.Range("M4").Value = 38.36920666666666
$ws.Range("N4").Value = 115.10762
$ws.Range("O4").Value = 0.3167137665810073
$ws.Range("P4").Value = 0.3167137665810074
$ws.Range("Q4").Value = 42.35589513668889
$ws.Range("R4").Value = 381.2030562302
$ws.Range("S4").Value = 0.00757934997928877
$ws.Range("T4").Value = 0.007579349979288772
$ws.Range("H5").Value = 70.73212899999999
$ws.Range("I5").Value = 0.5111277390235027
$ws.Range("J5").Value = 0.5111277390235027
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 13.10121233333333
$ws.Range("N5").Value = 39.303637
$ws.Range("O5").Value = 0.1081423012186565
$ws.Range("P5").Value = 0.1081423012186565
$ws.Range("Q5").Value = 308.8922136059081
$ws.Range("R5").Value = 2780.029922453173
$ws.Range("S5").Value = 0.05527452991469047
$ws.Range("T5").Value = 0.05527452991469047
$ws.Range("H6").Value = 70.73212899999999
$ws.Range("I6").Value = 0.5111277390235027
$ws.Range("J6").Value = 0.5111277390235027
$ws.Range("O6").Value = 0.5751439322003361
$ws.Range("P6").Value = 0.5751439322003362
$ws.Range("S6").Value = 0.2939720176786446
$ws.Range("T6").Value = 0.2939720176786446
$ws.Range("H7").Value = 70.73212899999999
$ws.Range("I7").Value = 0.5111277390235027
$ws.Range("J7").Value = 0.5111277390235027
$ws.Range("M7").Value = 38.36920666666666
$ws.Range("N7").Value = 115.10762
$ws.Range("O7").Value = 0.3167137665810073
$ws.Range("P7").Value = 0.3167137665810074
$ws.Range("Q7").Value = 904.645225191442
$ws.Range("R7").Value = 8141.807026722979
$ws.Range("S7").Value = 0.1618811914301677
$ws.Range("T7").Value = 0.1618811914301677
$ws.Range("G8").Value = 21.446869
$ws.Range("H8").Value = 64.34060699999999
$ws.Range("I8").Value = 0.4649410310173153
$ws.Range("J8").Value = 0.4649410310173154
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 13.10121233333333
$ws.Range("N8").Value = 39.303637
$ws.Range("O8").Value = 0.1081423012186565
$ws.Range("P8").Value = 0.1081423012186565
$ws.Range("Q8").Value = 280.9799846541843
$ws.Range("R8").Value = 2528.819861887659
$ws.Range("S8").Value = 0.05027979302518721
$ws.Range("T8").Value = 0.05027979302518722
$ws.Range("G9").Value = 21.446869
$ws.Range("H9").Value = 64.34060699999999
$ws.Range("I9").Value = 0.4649410310173153
$ws.Range("J9").Value = 0.4649410310173154
$ws.Range("O9").Value = 0.5751439322003361
$ws.Range("P9").Value = 0.5751439322003362
$ws.Range("Q9").Value = 1494.363735767425
$ws.Range("R9").Value = 13449.27362190683
$ws.Range("S9").Value = 0.2674080128205772
$ws.Range("T9").Value = 0.2674080128205772
$ws.Range("G10").Value = 21.446869
$ws.Range("H10").Value = 64.34060699999999
$ws.Range("I10").Value = 0.4649410310173153
$ws.Range("J10").Value = 0.4649410310173154
$ws.Range("M10").Value = 38.36920666666666
$ws.Range("N10").Value = 115.10762
$ws.Range("O10").Value = 0.3167137665810073
$ws.Range("P10").Value = 0.3167137665810074
$ws.Range("Q10").Value = 822.8993490139264
$ws.Range("R10").Value = 7406.094141125339
$ws.Range("S10").Value = 0.1472532251715509
$ws.Range("T10").Value = 0.147253225171551
